$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2: shared string content changes from "Desert Sage - 5 Item" to "3 MONTHS"
$ws.Range("A2").Value = "3 MONTHS"

# Update B2 / C2 numeric values
$ws.Range("B2").Value = 23729012754
$ws.Range("C2").Value = 177939546130

# Update D2: now references a new string "Print Paradise - 5 Items"
$ws.Range("D2").Value = "Print Paradise - 5 Items"

# Update selection to D2 to match sheetView selection change
$ws.Range("D2").Select()
